{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Summary of the change (see commit message \"Removed folder structure +\n// fixed bug with target groups\"):\n//  1. Remove the stray `_GoBack` bookmark that wrapped \"timestamp\" in the\n//     \"Eksportert: {{timestamp}}\" paragraph (leftover cursor-position\n//     bookmark, not meaningful content).\n//  2. Fix the `{{field.name}}` / `{{field.value}}` template placeholders in\n//     the table (row 2) so they read `{{ field.name }}` / `{{ field.value }}`\n//     (a padding-space bug fix so the templating engine resolves the\n//     \"target groups\" field correctly), splitting the runs and proofing\n//     marks (gramStart/gramEnd, spellStart/spellEnd) the way Word would\n//     after a manual retype, and leaving a fresh (empty) `_GoBack` bookmark\n//     behind - exactly like Word drops one wherever the cursor last was\n//     when the document was saved.\n\nconst body = context.document.body;\n\n// --- Step 1: drop the leftover `_GoBack` bookmark around \"timestamp\" -----\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// --- Step 2: locate the two placeholder paragraphs in the table ----------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet nameParagraph = null;\nlet valueParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text === \"{{field.name}}\") nameParagraph = p;\n  else if (p.text === \"{{field.value}}\") valueParagraph = p;\n}\n\nconst OOXML_NS =\n  'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction wrapPackage(bodyInnerXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    `<w:document ${OOXML_NS}><w:body>${bodyInnerXml}</w:body></w:document>` +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\nconst RPR = '<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:lang w:val=\"en-US\"/></w:rPr>';\nconst PPR = '<w:pPr><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>';\n\nfunction run(text, preserve) {\n  const sp = preserve ? ' xml:space=\"preserve\"' : \"\";\n  return `<w:r>${RPR}<w:t${sp}>${text}</w:t></w:r>`;\n}\n\nif (nameParagraph) {\n  const newPara =\n    \"<w:p>\" +\n    PPR +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    run(\"{{\") +\n    run(\" \", true) +\n    run(\"field.\") +\n    run(\"name\") +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    run(\" \", true) +\n    run(\"}}\") +\n    \"</w:p>\";\n  nameParagraph.insertOoxml(wrapPackage(newPara), \"Replace\");\n}\n\nif (valueParagraph) {\n  const newPara =\n    \"<w:p>\" +\n    PPR +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    run(\"{{\") +\n    run(\" \", true) +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    run(\"field\") +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    run(\".val\") +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    run(\"ue\") +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    run(\" \", true) +\n    run(\"}}\") +\n    \"</w:p>\";\n  valueParagraph.insertOoxml(wrapPackage(newPara), \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Summary of the change (see commit message \"Removed folder structure +\n# fixed bug with target groups\"):\n#  1. Remove the stray `_GoBack` bookmark that wrapped \"timestamp\" in the\n#     \"Eksportert: {{timestamp}}\" paragraph (leftover cursor-position\n#     bookmark, not meaningful content).\n#  2. Fix the `{{field.name}}` / `{{field.value}}` template placeholders in\n#     the table (row 2) so they read `{{ field.name }}` / `{{ field.value }}`\n#     (a padding-space bug fix so the templating engine resolves the\n#     \"target groups\" field correctly), splitting the runs and proofing\n#     marks (gramStart/gramEnd, spellStart/spellEnd) the way Word would\n#     after a manual retype, and leaving a fresh (empty) `_GoBack` bookmark\n#     behind - exactly like Word drops one wherever the cursor last was\n#     when the document was saved.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: drop the leftover `_GoBack` bookmark around \"timestamp\" ------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- Step 2: locate the two placeholder paragraphs in the table -----------\n$nameRange = $null\n$valueRange = $null\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    # Range.Text carries a trailing paragraph mark (chr 13) and, inside a\n    # table cell, a cell-mark (chr 7) too - strip both before comparing.\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"{{field.name}}\") {\n        $nameRange = $p.Range\n    } elseif ($t -eq \"{{field.value}}\") {\n        $valueRange = $p.Range\n    }\n}\n\nfunction New-FlatOpcPackage([string]$BodyInnerXml) {\n    $ns = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        \"<w:document $ns><w:body>$BodyInnerXml</w:body></w:document>\" +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\n$RPR = '<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:lang w:val=\"en-US\"/></w:rPr>'\n$PPR = '<w:pPr><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>'\n\nfunction New-Run([string]$Text, [bool]$Preserve = $false) {\n    if ($Preserve) {\n        return \"<w:r>$RPR<w:t xml:space=`\"preserve`\">$Text</w:t></w:r>\"\n    }\n    return \"<w:r>$RPR<w:t>$Text</w:t></w:r>\"\n}\n\nif ($nameRange -ne $null) {\n    $newPara = \"<w:p>\" + $PPR +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        (New-Run \"{{\") +\n        (New-Run \" \" $true) +\n        (New-Run \"field.\") +\n        (New-Run \"name\") +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n        (New-Run \" \" $true) +\n        (New-Run \"}}\") +\n        \"</w:p>\"\n    $nameRange.InsertXML((New-FlatOpcPackage $newPara))\n}\n\nif ($valueRange -ne $null) {\n    $newPara = \"<w:p>\" + $PPR +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        (New-Run \"{{\") +\n        (New-Run \" \" $true) +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        (New-Run \"field\") +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n        (New-Run \".val\") +\n        '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n        '<w:bookmarkEnd w:id=\"0\"/>' +\n        (New-Run \"ue\") +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        (New-Run \" \" $true) +\n        (New-Run \"}}\") +\n        \"</w:p>\"\n    $valueRange.InsertXML((New-FlatOpcPackage $newPara))\n}\n"}
